# Daily attendance processing - 2025-10-01 19:49:49
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: assign a value that *looks* numeric/date-like (percentages, dates,
# etc.) while forcing Excel to keep it as literal text, and re-apply the
# fill/font colors + centered alignment so the visual style matches the
# sheet's existing "stat" / "data" cell look even though Excel allocates a
# brand-new cell style index for the forced text format.
function Set-TextValue {
    param($range, [string]$text, $fillColor, $fontColor)
    $range.NumberFormat = "@"
    $range.Value = $text
    if ($fillColor -ne $null) { $range.Interior.Color = $fillColor }
    if ($fontColor -ne $null) { $range.Font.Color = $fontColor }
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4108
}

$yellowFill = 14745599   # matches fgColor 00FFFFE0 used by the K:S stats rows
$greenFill  = 9498256    # matches fgColor 0090EE90 used by the data rows
$blackFont  = 0          # matches font color 00000000

# --- Row 2 (B2A session 1) ---
$ws.Range("G2").Value = "backup@backdoor.com, system"
$ws.Range("H2").Value = "33/53"

# --- Row 6 (B2A session 5) ---
$ws.Range("G6").Value = "System, dnasr281@gmail.com"
$ws.Range("H6").Value = "44/53"
$ws.Range("L6").Value = 68

# --- Row 7 (Missing Sessions stat) ---
$ws.Range("L7").Value = 1

# --- Row 9 (Coverage %) ---
Set-TextValue $ws.Range("L9") "42.8%" $yellowFill $blackFont

# --- Row 10 (Average Attendance %) ---
Set-TextValue $ws.Range("L10") "62.2%" $yellowFill $blackFont

# --- Row 12 (B2A session 11) ---
$ws.Range("G12").Value = "System, dnasr281@gmail.com"
$ws.Range("H12").Value = "31/53"

# --- Row 13 (B2A session 12) ---
$ws.Range("G13").Value = "System, dnasr281@gmail.com"
$ws.Range("H13").Value = "36/53"

# --- Row 15 (Group stats B2A) ---
Set-TextValue $ws.Range("S15") "60.1%" $yellowFill $blackFont

# --- Row 18 (Group stats B2D) ---
$ws.Range("O18").Value = 11
$ws.Range("P18").Value = 0
Set-TextValue $ws.Range("R18") "42.3%" $yellowFill $blackFont
Set-TextValue $ws.Range("S18") "64.4%" $yellowFill $blackFont

# --- Row 19 (Group stats B2E) ---
$ws.Range("O19").Value = 11
$ws.Range("P19").Value = 0
Set-TextValue $ws.Range("R19") "42.3%" $yellowFill $blackFont
Set-TextValue $ws.Range("S19") "67.6%" $yellowFill $blackFont

# --- Row 20 (Group stats B2F) ---
$ws.Range("O20").Value = 11
$ws.Range("P20").Value = 0
Set-TextValue $ws.Range("R20") "42.3%" $yellowFill $blackFont
Set-TextValue $ws.Range("S20") "71.6%" $yellowFill $blackFont

# --- Row 29 (B2B session 1) ---
$ws.Range("G29").Value = "backup@backdoor.com, system"

# --- Row 33 (B2B session 5) ---
$ws.Range("G33").Value = "System, dnasr281@gmail.com"

# --- Row 39 (B2B session 11) ---
$ws.Range("G39").Value = "System, dnasr281@gmail.com"

# --- Row 40 (B2B session 12) ---
$ws.Range("G40").Value = "System, dnasr281@gmail.com"

# --- Row 56 (B2C session 1) ---
$ws.Range("G56").Value = "backup@backdoor.com, system"

# --- Row 60 (B2C session 5) ---
$ws.Range("G60").Value = "System, dnasr281@gmail.com"

# --- Row 66 (B2C session 11) ---
$ws.Range("G66").Value = "System, dnasr281@gmail.com"

# --- Row 67 (B2C session 12) ---
$ws.Range("G67").Value = "System, dnasr281@gmail.com"

# --- Row 90 (B2D session 8) ---
$ws.Range("G90").Value = "admin@admin.com, dnasr281@gmail.com"

# --- Row 92 (B2D session 10) ---
$ws.Range("H92").Value = "43/56"

# --- Row 93 (B2D session 11): was "Not Recorded" placeholder row, now filled in ---
# Copy formatting from row 92 (style s=2) onto row 93 (which currently has style s=9)
$ws.Range("A92:I92").Copy($ws.Range("A93:I93"))
$ws.Range("A93").Value = "Year 4"
$ws.Range("B93").Value = "B2D"
$ws.Range("C93").Value = "GENERAL SURGERY"
Set-TextValue $ws.Range("D93") "11" $greenFill $blackFont
Set-TextValue $ws.Range("E93") "01/10/2025" $greenFill $blackFont
$ws.Range("F93").Value = "10:30:00"
$ws.Range("G93").Value = "dnasr281@gmail.com"
$ws.Range("H93").Value = "46/56"
$ws.Range("I93").Value = "Recorded"

# --- Row 116 (B2E session 8) ---
$ws.Range("G116").Value = "admin@admin.com, dnasr281@gmail.com"

# --- Row 118 (B2E session 10) ---
$ws.Range("H118").Value = "45/55"

# --- Row 119 (B2E session 11): was "Not Recorded" placeholder row, now filled in ---
$ws.Range("A118:I118").Copy($ws.Range("A119:I119"))
$ws.Range("A119").Value = "Year 4"
$ws.Range("B119").Value = "B2E"
$ws.Range("C119").Value = "GENERAL SURGERY"
Set-TextValue $ws.Range("D119") "11" $greenFill $blackFont
Set-TextValue $ws.Range("E119") "01/10/2025" $greenFill $blackFont
$ws.Range("F119").Value = "10:30:00"
$ws.Range("G119").Value = "dnasr281@gmail.com"
$ws.Range("H119").Value = "40/55"
$ws.Range("I119").Value = "Recorded"

# --- Row 142 (B2F session 8) ---
$ws.Range("G142").Value = "admin@admin.com, dnasr281@gmail.com"

# --- Row 145 (B2F session 11): was "Not Recorded" placeholder row, now filled in ---
$ws.Range("A144:I144").Copy($ws.Range("A145:I145"))
$ws.Range("A145").Value = "Year 4"
$ws.Range("B145").Value = "B2F"
$ws.Range("C145").Value = "GENERAL SURGERY"
Set-TextValue $ws.Range("D145") "11" $greenFill $blackFont
Set-TextValue $ws.Range("E145") "01/10/2025" $greenFill $blackFont
$ws.Range("F145").Value = "10:30:00"
$ws.Range("G145").Value = "dnasr281@gmail.com"
$ws.Range("H145").Value = "48/57"
$ws.Range("I145").Value = "Recorded"
